$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.667443990707397
$ws.Range("B1").Value = 4.055765151977539
$ws.Range("C1").Value = 3.357496023178101
$ws.Range("D1").Value = 1.582585334777832
$ws.Range("E1").Value = 0.7899190187454224
